$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1302.5625
$ws.Range("I15").Value = 1302.5625
$ws.Range("K15").Value = 3907.6875
$ws.Range("M15").Value = -3738.6875
$ws.Range("H86").Value = 4993
$ws.Range("I86").Value = 4990
$ws.Range("J86").Value = 4994.5
$ws.Range("K86").Value = 4990
$ws.Range("L86").Value = 4994.5
$ws.Range("M86").Value = -3867
$ws.Range("N86").Value = -7240.5
$ws.Range("H89").Value = 4993
$ws.Range("I89").Value = 4990
$ws.Range("J89").Value = 4994.5
$ws.Range("K89").Value = 24950
$ws.Range("L89").Value = 24972.5
$ws.Range("M89").Value = -19334
$ws.Range("N89").Value = -36204.5
$ws.Range("H98").Value = 795.1177
$ws.Range("I98").Value = 913.2857
$ws.Range("J98").Value = 243.66667
$ws.Range("K98").Value = 913.2857
$ws.Range("L98").Value = 243.66667
$ws.Range("M98").Value = 584.7143
$ws.Range("N98").Value = -3239.66667
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").Value = ""
$ws.Range("H112").Value = 1817.3334
$ws.Range("J112").Value = 1989.7778
$ws.Range("L112").Value = 5969.3334
$ws.Range("N112").Value = -8185.3334
$ws.Range("H122").Value = 795.1177
$ws.Range("I122").Value = 913.2857
$ws.Range("J122").Value = 243.66667
$ws.Range("K122").Value = 2739.8571
$ws.Range("L122").Value = 731.00001
$ws.Range("M122").Value = -289.8571000000002
$ws.Range("N122").Value = -5631.00001
$ws.Range("H132").Value = 1152.4839
$ws.Range("I132").Value = 1176.4
$ws.Range("J132").Value = 435
$ws.Range("K132").Value = 3529.2
$ws.Range("L132").Value = 1305
$ws.Range("M132").Value = -999.2000000000003
$ws.Range("N132").Value = -6365
$ws.Range("H137").Value = 2477.5293
$ws.Range("I137").Value = 1494.95
$ws.Range("J137").Value = 3881.2144
$ws.Range("K137").Value = 4484.85
$ws.Range("L137").Value = 11643.6432
$ws.Range("M137").Value = -1934.85
$ws.Range("N137").Value = -16743.6432
$ws.Range("H138").Value = 4286.9116
$ws.Range("I138").Value = 3911.4375
$ws.Range("J138").Value = 4620.6665
$ws.Range("K138").Value = 11734.3125
$ws.Range("L138").Value = 13861.9995
$ws.Range("M138").Value = -6594.3125
$ws.Range("N138").Value = -24141.9995
$ws.Range("I141").Value = 7999.5
$ws.Range("J141").Value = 8499.5
$ws.Range("K141").Value = 23998.5
$ws.Range("L141").Value = 25498.5
$ws.Range("M141").Value = -18818.5
$ws.Range("N141").Value = -35858.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2233.375
$ws.Range("I2").Value = 2093.25
$ws.Range("K2").Value = 2093.25
$ws.Range("M2").Value = -1980.25
$ws.Range("H32").Value = 9215.309999999999
$ws.Range("I32").Value = 8870.179
$ws.Range("K32").Value = 8870.179
$ws.Range("M32").Value = -8583.179
$ws.Range("H61").Value = 1437.2354
$ws.Range("I61").Value = 1056.1111
$ws.Range("K61").Value = 1056.1111
$ws.Range("M61").Value = -844.1111000000001
$ws.Range("H74").Value = 1453.6389
$ws.Range("I74").Value = 1052.7931
$ws.Range("K74").Value = 1052.7931
$ws.Range("M74").Value = -178.7931000000001
$ws.Range("H77").Value = 1453.6389
$ws.Range("I77").Value = 1052.7931
$ws.Range("K77").Value = 5263.9655
$ws.Range("M77").Value = -895.9655000000002
$ws.Range("H116").Value = 2233.375
$ws.Range("I116").Value = 2093.25
$ws.Range("K116").Value = 2093.25
$ws.Range("M116").Value = 200.75
$ws.Range("H136").Value = 1437.2354
$ws.Range("I136").Value = 1056.1111
$ws.Range("K136").Value = 3168.3333
$ws.Range("M136").Value = -618.3333000000002
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2233.375
$ws.Range("I3").Value = 2093.25
$ws.Range("K3").Value = 2093.25
$ws.Range("M3").Value = -1979.25
$ws.Range("H20").Value = 10932.333
$ws.Range("I20").Value = 9899
$ws.Range("K20").Value = 9899
$ws.Range("M20").Value = -9652
$ws.Range("H105").Value = 6917.2
$ws.Range("I105").Value = 6917.2
$ws.Range("K105").Value = 6917.2
$ws.Range("M105").Value = -5170.2
$ws.Range("H107").Value = 1147
$ws.Range("J107").Value = 1999
$ws.Range("L107").Value = 1999
$ws.Range("N107").Value = -5839
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2970.3333
$ws.Range("I58").Value = 2948.5
$ws.Range("J58").Value = 3014
$ws.Range("K58").Value = 2948.5
$ws.Range("L58").Value = 3014
$ws.Range("M58").Value = -2745.5
$ws.Range("N58").Value = -3420
$ws.Range("H94").Value = 1110
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").Value = ""
$ws.Range("H132").Value = 4699
$ws.Range("I132").Value = 4499
$ws.Range("K132").Value = 13497
$ws.Range("M132").Value = -10967
$ws.Range("H136").Value = 2970.3333
$ws.Range("I136").Value = 2948.5
$ws.Range("J136").Value = 3014
$ws.Range("K136").Value = 8845.5
$ws.Range("L136").Value = 9042
$ws.Range("M136").Value = -6295.5
$ws.Range("N136").Value = -14142
$ws.Range("H141").Value = 51045.41
$ws.Range("J141").Value = 55184.8
$ws.Range("L141").Value = 55184.8
$ws.Range("N141").Value = -65544.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 101689.7
$ws.Range("J122").Value = 101689.7
$ws.Range("L122").Value = 915207.2999999999
$ws.Range("N122").Value = -920107.2999999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2912.25
$ws.Range("I102").Value = 2883.1667
$ws.Range("J102").Value = 2999.5
$ws.Range("K102").Value = 2883.1667
$ws.Range("L102").Value = 2999.5
$ws.Range("M102").Value = -1261.1667
$ws.Range("N102").Value = -6243.5
$ws.Range("H132").Value = 4673.3335
$ws.Range("I132").Value = 4651.643
$ws.Range("J132").Value = 4749.25
$ws.Range("K132").Value = 13954.929
$ws.Range("L132").Value = 14247.75
$ws.Range("M132").Value = -11424.929
$ws.Range("N132").Value = -19307.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5188.4443
$ws.Range("I22").Value = 3559.2
$ws.Range("J22").Value = 7225
$ws.Range("K22").Value = 3559.2
$ws.Range("L22").Value = 7225
$ws.Range("M22").Value = -3264.2
$ws.Range("N22").Value = -7815
$ws.Range("H27").Value = 5188.4443
$ws.Range("I27").Value = 3559.2
$ws.Range("J27").Value = 7225
$ws.Range("K27").Value = 3559.2
$ws.Range("L27").Value = 7225
$ws.Range("M27").Value = -3452.2
$ws.Range("N27").Value = -7439
$ws.Range("H68").Value = 2999
$ws.Range("I68").Value = 2999
$ws.Range("K68").Value = 2999
$ws.Range("M68").Value = -2250
$ws.Range("H71").Value = 2999
$ws.Range("I71").Value = 2999
$ws.Range("K71").Value = 14995
$ws.Range("M71").Value = -11251
$ws.Range("H136").Value = 4572.6523
$ws.Range("I136").Value = 4788.706
$ws.Range("J136").Value = 3960.5
$ws.Range("K136").Value = 14366.118
$ws.Range("L136").Value = 11881.5
$ws.Range("M136").Value = -11816.118
$ws.Range("N136").Value = -16981.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5184.9
$ws.Range("I81").Value = 3649.6667
$ws.Range("K81").Value = 7299.3334
$ws.Range("M81").Value = -6238.3334
$ws.Range("H84").Value = 5184.9
$ws.Range("I84").Value = 3649.6667
$ws.Range("K84").Value = 36496.667
$ws.Range("M84").Value = -31192.667
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = ""
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""
$ws.Range("H139").Value = 99000
$ws.Range("J139").Value = 99000
$ws.Range("L139").Value = 99000
$ws.Range("N139").Value = -109280
